$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Label" header in column H (mirror the formatting of the other header cells)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H1").Value = "Label"

# Populate column H with 0/1 "Label" values mirroring the boolean Success column (C)
# Block 1: rows 2-11 (Iterations = 100)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1

# Block 2: rows 12-21 (Iterations = 200)
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
